# Adds the new "extra objects" rows (naam/prioriteit) to the N:P mini-table
# on rows 14-19, matching the committed change ("dingen toegevoegd aan excel").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values (column N = naam, column P = prioriteit letter) ---------------
# (the source workbook pads these labels with leading spaces, kept verbatim)
$ws.Range("N14").Value = "             brug"
$ws.Range("N15").Value = "       lantaarnpaal"
$ws.Range("N16").Value = "     staande boom"
$ws.Range("N17").Value = "     grotere stenen"
$ws.Range("N18").Value = "        kampvuur"
$ws.Range("N19").Value = " stronken als stoelen"

$ws.Range("P14").Value = "S"
$ws.Range("P15").Value = "C"
$ws.Range("P16").Value = "C"
$ws.Range("P17").Value = "S"
$ws.Range("P18").Value = "C"
$ws.Range("P19").Value = "C"

# --- formatting -------------------------------------------------------------
# Priority letters get the same "should" (orange) / "could" (yellow) look
# already used throughout the D and P columns.
$ws.Range("D25").Copy()
$ws.Range("P14").PasteSpecial(-4122)

$ws.Range("D8").Copy()
$ws.Range("P15").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("P16").PasteSpecial(-4122)

$ws.Range("D25").Copy()
$ws.Range("P17").PasteSpecial(-4122)

$ws.Range("D8").Copy()
$ws.Range("P18").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("P19").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# keep the name column left unmerged (matches source) and move the
# selection to where the author last left off before saving.
$ws.Range("T17").Select()
